$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the (now unused) fill-applying style previously applied to D7:D12.
$ws.Range("D7:D12").ClearFormats()

# Add new test case row "sc16" (#1209 / oasislmf issue #1207).
$ws.Range("A20").Value() = "sc16"
$ws.Range("B20").Value() = 2
$ws.Range("C20").Value() = 2
$ws.Range("D20").Value() = 4
$ws.Range("E20").Value() = "A sublimit and a restriction on one of two policies"
$ws.Range("F20").Value() = "complete"
$ws.Range("G20").Value() = "yes"
$ws.Range("H20").Value() = "done"

# Update selection/view to the newly added row.
[void]$ws.Range("F20").Select()
